$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Year (C) and IssueDate (E) columns to stay text, not get auto-converted
# to numbers/dates, before writing the new values into them.
$ws.Range("C2:C6").NumberFormat = "@"
$ws.Range("E2:E6").NumberFormat = "@"

# --- Update data rows 2-6 with new November SEBI data ---
# Row 2
$ws.Range("A2").Value = 'Listed Companies'
$ws.Range("B2").Value = 'Circular-BSE'
$ws.Range("C2").Value = '2025'
$ws.Range("D2").Value = 'November'
$ws.Range("E2").Value = '2025-11-21'
$ws.Range("F2").Value = 'Master Circular for issue and listing of Non-convertible Securities, Securitised Debt Instruments, Security Receipts, Municipal Debt Securities and Commercial Paper.'
$ws.Range("G2").Value = 'https://www.bseindia.com/markets/MarketInfo/DownloadAttach.aspx?id=20251121-73&attachedId=80bb176a-9beb-47e3-99ff-4832d65a2c1b'
$ws.Range("H2").Value = 'Master_Circular_for_issue_and_listing_of_Non_convertible_Securities_Securitised_Debt_Instruments_Sec.pdf'
$ws.Range("I2").Value = '/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/Listed Companies/Circular-BSE/2025/November/Master_Circular_for_issue_and_listing_of_Non_convertible_Securities_Securitised_Debt_Instruments_Sec.pdf'

# Row 3
$ws.Range("A3").Value = 'Listed Companies'
$ws.Range("B3").Value = 'Circular-BSE'
$ws.Range("C3").Value = '2025'
$ws.Range("D3").Value = 'November'
$ws.Range("E3").Value = '2025-11-21'
$ws.Range("F3").Value = 'XBRL based filing of Regulation 50 for Entities which has listed its non-convertible securities'
$ws.Range("G3").Value = 'PrintToPDF'
$ws.Range("H3").Value = 'XBRL_based_filing_of_Regulation_50_for_Entities_which_has_listed_its_non_convertible_securities.pdf'
$ws.Range("I3").Value = '/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/Listed Companies/Circular-BSE/2025/November/XBRL_based_filing_of_Regulation_50_for_Entities_which_has_listed_its_non_convertible_securities.pdf'

# Row 4
$ws.Range("A4").Value = 'Listed Companies'
$ws.Range("B4").Value = 'Circular-BSE'
$ws.Range("C4").Value = '2025'
$ws.Range("D4").Value = 'November'
$ws.Range("E4").Value = '2025-11-17'
$ws.Range("F4").Value = 'Frequently Asked Questions (FAQ) for submission of financial results as required under Regulation 33 of SEBI (LODR) Regulations, 2015 & Master circular for compliance with the provisions of the SEBI (LODR) Regulations, 2015.'
$ws.Range("G4").Value = 'https://www.bseindia.com/markets/MarketInfo/DownloadAttach.aspx?id=20251117-20&attachedId=7f776d85-62d5-4358-8d03-694f1de8f401'
$ws.Range("H4").Value = 'Frequently_Asked_Questions_FAQ_for_submission_of_financial_results_as_required_under_Regulation_33_o.pdf'
$ws.Range("I4").Value = '/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/Listed Companies/Circular-BSE/2025/November/Frequently_Asked_Questions_FAQ_for_submission_of_financial_results_as_required_under_Regulation_33_o.pdf'

# Row 5
$ws.Range("A5").Value = 'SEBI'
$ws.Range("B5").Value = 'Press Release'
$ws.Range("C5").Value = '2025'
$ws.Range("D5").Value = 'November'
$ws.Range("E5").Value = '2025-11-19'
$ws.Range("F5").Value = 'Caution to Public regarding unregistered Online Bond Platform Providers'
$ws.Range("G5").Value = 'https://www.sebi.gov.in/sebi_data/attachdocs/nov-2025/1763551749033.pdf'
$ws.Range("H5").Value = '1763551749033.pdf'
$ws.Range("I5").Value = '/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/SEBI/Press Release/2025/November/1763551749033.pdf'

# Row 6
$ws.Range("A6").Value = 'SEBI'
$ws.Range("B6").Value = 'Consulatation Paper'
$ws.Range("C6").Value = '2025'
$ws.Range("D6").Value = 'November'
$ws.Range("E6").Value = '2025-11-19'
$ws.Range("F6").Value = 'Extension of timeline for submission of public comments on the consultation paper on ''Comprehensive review of SEBI (Mutual Funds) regulations, 1996'''
$ws.Range("G6").Value = 'https://www.sebi.gov.in/sebi_data/attachdocs/nov-2025/1763551731559.pdf'
$ws.Range("H6").Value = '1763551731559.pdf'
$ws.Range("I6").Value = '/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/SEBI/Consulatation Paper/2025/November/1763551731559.pdf'

# Restore plain "Normal" style on the Year/IssueDate columns (removes the
# temporary text-number-format so cells end up styleless, like the rest).
$ws.Range("C2:C6").Style = "Normal"
$ws.Range("E2:E6").Style = "Normal"

# --- Remove extra row 7 (now unused - only 5 data rows remain) ---
$ws.Rows("7").Delete()

# --- G3 (XBRL filing) has no real hyperlink; strip any leftover hyperlink styling ---
$ws.Range("G3").Style = "Normal"

# --- Rebuild hyperlinks collection to match new targets (G2, G4, G5, G6) ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("G2"), 'https://www.bseindia.com/markets/MarketInfo/DownloadAttach.aspx?id=20251121-73&attachedId=80bb176a-9beb-47e3-99ff-4832d65a2c1b')
$ws.Hyperlinks.Add($ws.Range("G4"), 'https://www.bseindia.com/markets/MarketInfo/DownloadAttach.aspx?id=20251117-20&attachedId=7f776d85-62d5-4358-8d03-694f1de8f401')
$ws.Hyperlinks.Add($ws.Range("G5"), 'https://www.sebi.gov.in/sebi_data/attachdocs/nov-2025/1763551749033.pdf')
$ws.Hyperlinks.Add($ws.Range("G6"), 'https://www.sebi.gov.in/sebi_data/attachdocs/nov-2025/1763551731559.pdf')
